$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert two new columns (B and C) for "Title" and "Description" ---
# Existing B..H shift right to D..J automatically; data validations on the
# plain (non-extLst) list shift with them.
$ws.Range("B1:C1").EntireColumn.Insert()

# --- Header row ---
$ws.Cells.Item(1,2).Value = "Title"
$ws.Cells.Item(1,3).Value = "Description"

# --- Row 2: patient name question ---
$ws.Cells.Item(2,2).Value = "Patient"
$ws.Cells.Item(2,3).Value = "Please enter Name…."

# --- Row 3: patient age question ---
$ws.Cells.Item(3,2).Value = "Age"
$ws.Cells.Item(3,3).Value = "Please enter Age …."

# --- Row 4: patient sex question ---
$ws.Cells.Item(4,2).Value = "Sex"
$ws.Cells.Item(4,3).Value = "Please enter Sex…."

# --- Row 5: pregnancy question ---
$ws.Cells.Item(5,2).Value = "Pregnant"
$ws.Cells.Item(5,3).Value = "Please enter Pregnant..."

# --- Row 6: symptoms / onset date question ---
$ws.Cells.Item(6,2).Value = "symptoms"
$ws.Cells.Item(6,3).Value = "Please enter symptoms…"

# --- Rename the If_Condition branch headers ---
$ws.Cells.Item(1,9).Value = "Then_Goto"
$ws.Cells.Item(1,10).Value = "Else_Goto"

# NOTE: the boolean-required list validation on column D auto-shifts to F
# (and the DataTypes!$A:$A extended validation remains tied to its sqref)
# when the columns are inserted above, so nothing further is needed here.

# --- Column widths for the newly inserted Title/Description columns ---
$ws.Range("B1:C1").ColumnWidth = 26.109375

# --- Selection moved to I1 ---
$null = $ws.Range("I1").Select()
